# Apply Introducer name updates / renames and one commission correction
# to the "Introducer Commissions" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Introducer (column B) renames ---

# "Anand Sethia" -> "Setcap"
$ws.Range("B3:B4").Value = "Setcap"
$ws.Range("B283:B284").Value = "Setcap"

# "Rick + Andrew" -> "Altras Capital Financing Broker"
$ws.Range("B332:B333").Value = "Altras Capital Financing Broker"

# "Anand" -> "Setcap"
$ws.Range("B378:B392").Value = "Setcap"
$ws.Range("B505:B519").Value = "Setcap"

# "Dan" -> "Daniel Baumslag"
$ws.Range("B393:B411").Value = "Daniel Baumslag"
$ws.Range("B467:B478").Value = "Daniel Baumslag"
$ws.Range("B487:B493").Value = "Daniel Baumslag"

# "Anand+Dan" -> "Setcap+Daniel Baumslag"
$ws.Range("B520:B521").Value = "Setcap+Daniel Baumslag"

# --- Commission correction on row 562 (VC133 / Altras Capital Financing Broker) ---
# Rate (bps): 25 -> 500
$ws.Range("E562").Value = 500

# Rate (%) and Commission Amount are stored as text in this sheet, so force
# the cells to Text format before assigning to avoid Excel re-interpreting
# the values as numbers (which would drop the fixed decimal formatting).
$ws.Range("F562").NumberFormat = "@"
$ws.Range("F562").Value = "5.00"

$ws.Range("G562").NumberFormat = "@"
$ws.Range("G562").Value = "50000.00"
